$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 14
$ws.Range("F4").Value = 95
$ws.Range("F5").Value = 914
$ws.Range("F6").Value = 61
$ws.Range("F7").Value = 7156
$ws.Range("F8").Value = 112
$ws.Range("F9").Value = 158
$ws.Range("F10").Value = 6556
$ws.Range("F11").Value = 133
$ws.Range("F13").Value = 4533
$ws.Range("F17").Value = 4591
$ws.Range("F18").Value = 15
$ws.Range("F21").Value = 350
$ws.Range("F28").Value = 8210
$ws.Range("F30").Value = 1425
$ws.Range("F31").Value = 60
$ws.Range("F32").Value = 725
$ws.Range("F34").Value = 51
$ws.Range("F35").Value = 989
$ws.Range("F37").Value = 1685
$ws.Range("F39").Value = 968
$ws.Range("F41").Value = 4253
$ws.Range("F43").Value = 630
$ws.Range("F44").Value = 119
$ws.Range("F46").Value = 851
$ws.Range("F47").Value = 1131
$ws.Range("F48").Value = 7
$ws.Range("F49").Value = 25

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 24

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 239

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 239
$ws.Range("F3").Value = 24
$ws.Range("F5").Value = 14
$ws.Range("F7").Value = 95
$ws.Range("F9").Value = 914
$ws.Range("F10").Value = 61
$ws.Range("F11").Value = 7156
$ws.Range("F12").Value = 112
$ws.Range("F13").Value = 158
$ws.Range("F14").Value = 6556
$ws.Range("F15").Value = 133
$ws.Range("F17").Value = 4533
$ws.Range("F21").Value = 4591
$ws.Range("F22").Value = 15
$ws.Range("F24").Value = 350
$ws.Range("F29").Value = 8210
$ws.Range("F31").Value = 1425
$ws.Range("F32").Value = 60
$ws.Range("F33").Value = 725
$ws.Range("F35").Value = 51
$ws.Range("F36").Value = 989
$ws.Range("F37").Value = 1685
$ws.Range("F39").Value = 968
$ws.Range("F41").Value = 4253
$ws.Range("F43").Value = 630
$ws.Range("F44").Value = 119
$ws.Range("F46").Value = 851
$ws.Range("F47").Value = 1131
$ws.Range("F48").Value = 7
$ws.Range("F49").Value = 25
